$wb = $excel.ActiveWorkbook

# Fix the "PostId" DataType typo on the Comments sheet: "interger" -> "integer"
$wsComments = $wb.Worksheets.Item("Comments")
$wsComments.Cells.Item(4, 2).Value = "integer"

# Make the Comments sheet the active tab, with B5 selected
$wsComments.Activate()
[void]$wsComments.Range("B5").Select()
